# Update "想去人数" (want-to-go count) figures that changed between scrapes.
# These edits mirror a re-scrape of the source site: the counts for a
# given event were bumped (in "展览" / "本地生活" / "全部类型" sheets) while
# "演出" is untouched.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 8473
$ws1.Range("F4").Value  = 1940
$ws1.Range("F5").Value  = 6539
$ws1.Range("F7").Value  = 2100
$ws1.Range("F8").Value  = 585
$ws1.Range("F11").Value = 60
$ws1.Range("F16").Value = 8654
$ws1.Range("F21").Value = 1818
$ws1.Range("F26").Value = 56
$ws1.Range("F28").Value = 192
$ws1.Range("F29").Value = 1012
$ws1.Range("F30").Value = 22
$ws1.Range("F31").Value = 20
$ws1.Range("F32").Value = 15
$ws1.Range("F33").Value = 2147
$ws1.Range("F34").Value = 859
$ws1.Range("F35").Value = 501
$ws1.Range("F39").Value = 217
$ws1.Range("F42").Value = 62
$ws1.Range("F45").Value = 3977

# --- Sheet "本地生活" (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 2333
$ws3.Range("F4").Value = 315

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 2333
$ws4.Range("F6").Value  = 8473
$ws4.Range("F8").Value  = 315
$ws4.Range("F9").Value  = 1940
$ws4.Range("F10").Value = 6539
$ws4.Range("F11").Value = 2100
$ws4.Range("F13").Value = 585
$ws4.Range("F17").Value = 60
$ws4.Range("F20").Value = 8654
$ws4.Range("F24").Value = 1818
$ws4.Range("F28").Value = 56
$ws4.Range("F30").Value = 192
$ws4.Range("F31").Value = 22
$ws4.Range("F32").Value = 20
$ws4.Range("F33").Value = 15
$ws4.Range("F34").Value = 2147
$ws4.Range("F35").Value = 859
$ws4.Range("F37").Value = 501
$ws4.Range("F40").Value = 217
$ws4.Range("F44").Value = 3977

$wb.Save()
